# Insert a new weekly record at row 237 for "Hortaliza, Macroferia Regional de
# Talca - Repollo" (the existing rows 237-298 all shift down to 238-299).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 237:298 down by one row, duplicating formatting from the row above.
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new weekly observation.
$ws.Range("A237").Value = 5
$ws.Range("B237").Value = "Macroferia Regional de Talca"
$ws.Range("C237").Value = "Maule"
$ws.Range("D237").Value = 44736
$ws.Range("E237").Value = 7
$ws.Range("F237").Value = 100112006
$ws.Range("G237").Value = "Repollo"
$ws.Range("H237").Value = "Crespo record"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 3000
$ws.Range("K237").Value = 1100
$ws.Range("L237").Value = 1100
$ws.Range("M237").Value = 1100
$ws.Range("N237").Value = "$/unidad"
$ws.Range("O237").Value = "Región del Maule"
$ws.Range("P237").Value = 1100
$ws.Range("Q237").Value = 1
$ws.Range("R237").Value = "Hortaliza"
